# Update forecast values in the active worksheet (Sheet1).
# These cells hold model output ("Random Forest", column C) and actual
# observed values ("Gerçek Değerler", column D) keyed by Day/TIMEPRD.
# The new upload refreshes a handful of previously-zero (or stale) rows
# with recomputed figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C10" = 1
    "D10" = 1
    "C11" = 1
    "D11" = 1

    "C33" = 0.01666666666666667

    "C34" = 1
    "D34" = 1
    "C35" = 1
    "D35" = 1
    "C36" = 2.054545454545455

    "C58" = 1
    "D58" = 1
    "C59" = 0.9666666666666667
    "D59" = 1
    "C60" = 2.357142857142857

    "C82" = 1
    "D82" = 1
    "C83" = 1
    "D83" = 1

    "C106" = 1
    "D106" = 1
    "C107" = 1
    "D107" = 1

    "C130" = 1
    "D130" = 1
    "C131" = 1
    "D131" = 1

    "C154" = 1
    "D154" = 1
    "C155" = 1
    "D155" = 1
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
